$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new member rows after the existing data (rows 12 and 13)
$ws.Range("A12").Value = "Qananiisaa Biqilaa"
$ws.Range("B12").Value = 911818685
$ws.Range("C12").Value = 10000
$ws.Range("D12").Value = 100000

$ws.Range("A13").Value = "Oromiyaa Walfaanaa"
$ws.Range("B13").Value = 912861288
$ws.Range("C13").Value = 10000
$ws.Range("D13").Value = 100000

# Update the active selection to reflect where the user ended up (D14)
$ws.Range("D14").Select()
